$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Week 6 hours - row 20 (Formal team meeting)
$ws.Range("A20").Value = 6
$ws.Range("B20").Value2 = 43704
$ws.Range("C20").Value2 = 0.39583333333333331
$ws.Range("D20").Value2 = 43704
$ws.Range("E20").Value2 = 0.40972222222222227
$ws.Range("F20").Value = "Formal team meeting"

# Row 22 text entered before row 21's so the shared-string table gets the
# same ordering as the authored workbook (idx 23 = row22 text, idx 24 = row21 text)
$ws.Range("F22").Value = "Creating user database + authentication "
$ws.Range("F21").Value = "Began registration form "

# Week 6 hours - row 21 (Began registration form)
$ws.Range("A21").Value = 6
$ws.Range("B21").Value2 = 43706
$ws.Range("C21").Value2 = 0.83333333333333337
$ws.Range("D21").Value2 = 43706
$ws.Range("E21").Value2 = 0.875

# Week 6 hours - row 22 (Creating user database + authentication)
$ws.Range("A22").Value = 6
$ws.Range("B22").Value2 = 43707
$ws.Range("C22").Value2 = 0.54166666666666663
$ws.Range("D22").Value2 = 43707
$ws.Range("E22").Value2 = 0.625
